# Scheduled runner update: refresh market-derived profit figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H-N)
# across the per-job "Ultima Profits" sheets. Only pre-computed numeric
# values are touched - no formulas, headers, or formatting are affected.

$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H64").Value = 6996234.5
$ws.Range("I64").Value = 10991869
$ws.Range("J64").Value = 3875
$ws.Range("K64").Value = 10991869
$ws.Range("L64").Value = 3875
$ws.Range("M64").Value = -10991621
$ws.Range("N64").Value = -4371

$ws.Range("H67").Value = 6996234.5
$ws.Range("I67").Value = 10991869
$ws.Range("J67").Value = 3875
$ws.Range("K67").Value = 10991869
$ws.Range("L67").Value = 3875
$ws.Range("M67").Value = -10991011
$ws.Range("N67").Value = -5591

$ws.Range("H70").Value = 1678.5714
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 1692.3077
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 5076.9231
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -5616.9231

$ws.Range("H73").Value = 1678.5714
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 1692.3077
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 5076.9231
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -6948.9231

$ws.Range("H76").Value = 5413.5
$ws.Range("I76").Value = 3500
$ws.Range("K76").Value = 3500
$ws.Range("M76").Value = -3185

$ws.Range("H79").Value = 5413.5
$ws.Range("I79").Value = 3500
$ws.Range("K79").Value = 3500
$ws.Range("M79").Value = -2408

$ws.Range("H137").Value = 5263807.5
$ws.Range("I137").Value = 536.23334
$ws.Range("K137").Value = 1608.70002
$ws.Range("M137").Value = 941.29998

$ws.Range("H141").Value = 1101.2554
$ws.Range("I141").Value = 1059.9783
$ws.Range("K141").Value = 3179.9349
$ws.Range("M141").Value = 2000.0651

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 8197855.5
$ws.Range("I74").Value = 10639364
$ws.Range("J74").Value = 1364.4286
$ws.Range("K74").Value = 10639364
$ws.Range("L74").Value = 1364.4286
$ws.Range("M74").Value = -10638490
$ws.Range("N74").Value = -3112.4286

$ws.Range("H77").Value = 8197855.5
$ws.Range("I77").Value = 10639364
$ws.Range("J77").Value = 1364.4286
$ws.Range("K77").Value = 53196820
$ws.Range("L77").Value = 6822.143
$ws.Range("M77").Value = -53192452
$ws.Range("N77").Value = -15558.143

$ws.Range("H88").Value = 3084.8333
$ws.Range("I88").Value = 3057.1428
$ws.Range("J88").Value = 3102.4546
$ws.Range("K88").Value = 3057.1428
$ws.Range("L88").Value = 3102.4546
$ws.Range("M88").Value = -2651.1428
$ws.Range("N88").Value = -3914.4546

$ws.Range("H91").Value = 3084.8333
$ws.Range("I91").Value = 3057.1428
$ws.Range("J91").Value = 3102.4546
$ws.Range("K91").Value = 3057.1428
$ws.Range("L91").Value = 3102.4546
$ws.Range("M91").Value = -1653.1428
$ws.Range("N91").Value = -5910.4546

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 29413950
$ws.Range("I86").Value = 1833.3334
$ws.Range("J86").Value = 100003030
$ws.Range("K86").Value = 1833.3334
$ws.Range("L86").Value = 100003030
$ws.Range("M86").Value = -710.3334
$ws.Range("N86").Value = -100005276

$ws.Range("H89").Value = 29413950
$ws.Range("I89").Value = 1833.3334
$ws.Range("J89").Value = 100003030
$ws.Range("K89").Value = 9166.666999999999
$ws.Range("L89").Value = 500015150
$ws.Range("M89").Value = -3550.666999999999
$ws.Range("N89").Value = -500026382

$ws.Range("H105").Value = 4840.5625
$ws.Range("I105").Value = 3724.5
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3724.5
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1977.5
$ws.Range("N105").Value = -8494

$ws.Range("H134").Value = 2908.7556
$ws.Range("I134").Value = 2499.487
$ws.Range("J134").Value = 5569
$ws.Range("K134").Value = 7498.461
$ws.Range("L134").Value = 16707
$ws.Range("M134").Value = -4963.461
$ws.Range("N134").Value = -21777

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H62").Value = 2352
$ws.Range("I62").Value = 2200
$ws.Range("K62").Value = 2200
$ws.Range("M62").Value = -1576

$ws.Range("H65").Value = 2352
$ws.Range("I65").Value = 2200
$ws.Range("K65").Value = 11000
$ws.Range("M65").Value = -7880

$ws.Range("H94").Value = 3776.6667
$ws.Range("I94").Value = 2310.182
$ws.Range("J94").Value = 5017.5386
$ws.Range("K94").Value = 2310.182
$ws.Range("L94").Value = 5017.5386
$ws.Range("M94").Value = -1859.182
$ws.Range("N94").Value = -5919.5386

$ws.Range("H132").Value = 8198225
$ws.Range("I132").Value = 11112425
$ws.Range("J132").Value = 2037.0625
$ws.Range("K132").Value = 33337275
$ws.Range("L132").Value = 6111.1875
$ws.Range("M132").Value = -33334745
$ws.Range("N132").Value = -11171.1875

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H131").Value = 816.79
$ws.Range("I131").Value = 383.33334
$ws.Range("J131").Value = 859.65936
$ws.Range("K131").Value = 1150.00002
$ws.Range("L131").Value = 2578.97808
$ws.Range("M131").Value = 3889.99998
$ws.Range("N131").Value = -12658.97808

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 17226
$ws.Range("I70").Value = 19153.076
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 19153.076
$ws.Range("L70").Value = 4700
$ws.Range("M70").Value = -18883.076
$ws.Range("N70").Value = -5240

$ws.Range("H73").Value = 17226
$ws.Range("I73").Value = 19153.076
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 19153.076
$ws.Range("L73").Value = 4700
$ws.Range("M73").Value = -18217.076
$ws.Range("N73").Value = -6572

$ws.Range("H98").Value = 23821.5
$ws.Range("J98").Value = 23821.5
$ws.Range("L98").Value = 23821.5
$ws.Range("N98").Value = -29811.5

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H81").Value = 35611.668
$ws.Range("J81").Value = 35611.668
$ws.Range("L81").Value = 35611.668
$ws.Range("N81").Value = -37607.668

$ws.Range("H84").Value = 35611.668
$ws.Range("J84").Value = 35611.668
$ws.Range("L84").Value = 106835.004
$ws.Range("N84").Value = -116819.004

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 16086
$ws.Range("J62").Value = 22925.75
$ws.Range("L62").Value = 22925.75
$ws.Range("N62").Value = -24173.75

$ws.Range("H65").Value = 16086
$ws.Range("J65").Value = 22925.75
$ws.Range("L65").Value = 114628.75
$ws.Range("N65").Value = -120868.75

Write-Host "Applied 163 cell updates across 8 sheets / 28 rows"
